$wb = $excel.ActiveWorkbook

# --- Helper: locate source sheets for header-style copy ---
$pointsWs = $wb.Worksheets.Item("Points")
$avgPointsWs = $wb.Worksheets.Item("Avg Points")

# --- Insert "Rebounds" sheet right after "Assists" ---
$assistsWs = $wb.Worksheets.Item("Assists")
$reboundsWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $assistsWs)
$reboundsWs.Name = "Rebounds"

# --- Insert "3PM" sheet right after "Rebounds" ---
$pm3Ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $reboundsWs)
$pm3Ws.Name = "3PM"

# --- Insert "Avg Rebounds" sheet right after "Avg Assists" ---
$avgAssistsWs = $wb.Worksheets.Item("Avg Assists")
$avgReboundsWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $avgAssistsWs)
$avgReboundsWs.Name = "Avg Rebounds"

# --- Insert "Avg 3PM" sheet right after "Avg Rebounds" ---
$avgPm3Ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $avgReboundsWs)
$avgPm3Ws.Name = "Avg 3PM"

# --- Copy header row (with bold/border/style) from "Points" into "Rebounds" and "3PM" ---
$pointsWs.Range("A1:P1").Copy($reboundsWs.Range("A1:P1"))
$pointsWs.Range("A1:P1").Copy($pm3Ws.Range("A1:P1"))

# --- Copy header row (with bold/border/style) from "Avg Points" into "Avg Rebounds" and "Avg 3PM" ---
$avgPointsWs.Range("A1:B1").Copy($avgReboundsWs.Range("A1:B1"))
$avgPointsWs.Range("A1:B1").Copy($avgPm3Ws.Range("A1:B1"))
$avgReboundsWs.Range("B1").Value = "Avg Rebounds"
$avgPm3Ws.Range("B1").Value = "Avg 3PM"

# --- Fill "Rebounds" data rows 2-12 ---
$reboundsWs.Range("A2").Value = "'2025-10-21"
$reboundsWs.Range("B2").Value = "LAL"
$reboundsWs.Range("C2").Value = 2
$reboundsWs.Range("D2").Value = 9
$reboundsWs.Range("E2").Value = 7
$reboundsWs.Range("F2").Value = 1
$reboundsWs.Range("G2").Value = 0
$reboundsWs.Range("H2").Value = 1
$reboundsWs.Range("I2").Value = 5
$reboundsWs.Range("J2").Value = 0
$reboundsWs.Range("K2").Value = 5
$reboundsWs.Range("L2").Value = 2
$reboundsWs.Range("M2").Value = 7
$reboundsWs.Range("N2").Value = 1
$reboundsWs.Range("O2").Value = 0
$reboundsWs.Range("P2").Value = 0
$reboundsWs.Range("A3").Value = "'2025-10-23"
$reboundsWs.Range("B3").Value = "DEN"
$reboundsWs.Range("C3").Value = 6
$reboundsWs.Range("D3").Value = 5
$reboundsWs.Range("E3").Value = 3
$reboundsWs.Range("F3").Value = 5
$reboundsWs.Range("G3").Value = 0
$reboundsWs.Range("H3").Value = 2
$reboundsWs.Range("I3").Value = 5
$reboundsWs.Range("J3").Value = 0
$reboundsWs.Range("K3").Value = 3
$reboundsWs.Range("L3").Value = 0
$reboundsWs.Range("M3").Value = 8
$reboundsWs.Range("N3").Value = 6
$reboundsWs.Range("O3").Value = 0
$reboundsWs.Range("P3").Value = 0
$reboundsWs.Range("A4").Value = "'2025-10-24"
$reboundsWs.Range("B4").Value = "POR"
$reboundsWs.Range("C4").Value = 1
$reboundsWs.Range("D4").Value = 8
$reboundsWs.Range("E4").Value = 5
$reboundsWs.Range("F4").Value = 2
$reboundsWs.Range("G4").Value = 2
$reboundsWs.Range("H4").Value = 0
$reboundsWs.Range("I4").Value = 3
$reboundsWs.Range("J4").Value = 1
$reboundsWs.Range("K4").Value = 0
$reboundsWs.Range("L4").Value = 4
$reboundsWs.Range("M4").Value = 3
$reboundsWs.Range("N4").Value = 6
$reboundsWs.Range("O4").Value = 4
$reboundsWs.Range("P4").Value = 2
$reboundsWs.Range("A5").Value = "'2025-10-27"
$reboundsWs.Range("B5").Value = "MEM"
$reboundsWs.Range("C5").Value = 0
$reboundsWs.Range("D5").Value = 10
$reboundsWs.Range("E5").Value = 4
$reboundsWs.Range("F5").Value = 2
$reboundsWs.Range("G5").Value = 2
$reboundsWs.Range("H5").Value = 2
$reboundsWs.Range("I5").Value = 5
$reboundsWs.Range("J5").Value = 1
$reboundsWs.Range("K5").Value = 0
$reboundsWs.Range("L5").Value = 4
$reboundsWs.Range("M5").Value = 4
$reboundsWs.Range("N5").Value = 4
$reboundsWs.Range("O5").Value = 1
$reboundsWs.Range("P5").Value = 0
$reboundsWs.Range("A6").Value = "'2025-10-28"
$reboundsWs.Range("B6").Value = "LAC"
$reboundsWs.Range("C6").Value = 0
$reboundsWs.Range("D6").Value = 5
$reboundsWs.Range("E6").Value = 5
$reboundsWs.Range("F6").Value = 0
$reboundsWs.Range("G6").Value = 3
$reboundsWs.Range("H6").Value = 3
$reboundsWs.Range("I6").Value = 5
$reboundsWs.Range("J6").Value = 0
$reboundsWs.Range("K6").Value = 4
$reboundsWs.Range("L6").Value = 8
$reboundsWs.Range("M6").Value = 7
$reboundsWs.Range("N6").Value = 2
$reboundsWs.Range("O6").Value = 0
$reboundsWs.Range("P6").Value = 0
$reboundsWs.Range("A7").Value = "'2025-10-30"
$reboundsWs.Range("B7").Value = "MIL"
$reboundsWs.Range("C7").Value = 2
$reboundsWs.Range("D7").Value = 8
$reboundsWs.Range("E7").Value = 5
$reboundsWs.Range("F7").Value = 0
$reboundsWs.Range("G7").Value = 3
$reboundsWs.Range("H7").Value = 1
$reboundsWs.Range("I7").Value = 10
$reboundsWs.Range("J7").Value = 0
$reboundsWs.Range("K7").Value = 2
$reboundsWs.Range("L7").Value = 3
$reboundsWs.Range("M7").Value = 4
$reboundsWs.Range("N7").Value = 6
$reboundsWs.Range("O7").Value = 0
$reboundsWs.Range("P7").Value = 0
$reboundsWs.Range("A8").Value = "'2025-11-01"
$reboundsWs.Range("B8").Value = "IND"
$reboundsWs.Range("C8").Value = 3
$reboundsWs.Range("D8").Value = 5
$reboundsWs.Range("E8").Value = 4
$reboundsWs.Range("F8").Value = 0
$reboundsWs.Range("G8").Value = 2
$reboundsWs.Range("H8").Value = 1
$reboundsWs.Range("I8").Value = 6
$reboundsWs.Range("J8").Value = 4
$reboundsWs.Range("K8").Value = 6
$reboundsWs.Range("L8").Value = 3
$reboundsWs.Range("M8").Value = 10
$reboundsWs.Range("N8").Value = 0
$reboundsWs.Range("O8").Value = 0
$reboundsWs.Range("P8").Value = 0
$reboundsWs.Range("A9").Value = "'2025-11-04"
$reboundsWs.Range("B9").Value = "PHX"
$reboundsWs.Range("C9").Value = 0
$reboundsWs.Range("D9").Value = 5
$reboundsWs.Range("E9").Value = 2
$reboundsWs.Range("F9").Value = 0
$reboundsWs.Range("G9").Value = 5
$reboundsWs.Range("H9").Value = 2
$reboundsWs.Range("I9").Value = 4
$reboundsWs.Range("J9").Value = 0
$reboundsWs.Range("K9").Value = 0
$reboundsWs.Range("L9").Value = 7
$reboundsWs.Range("M9").Value = 5
$reboundsWs.Range("N9").Value = 4
$reboundsWs.Range("O9").Value = 2
$reboundsWs.Range("P9").Value = 3
$reboundsWs.Range("A10").Value = "'2025-11-05"
$reboundsWs.Range("B10").Value = "SAC"
$reboundsWs.Range("C10").Value = 2
$reboundsWs.Range("D10").Value = 9
$reboundsWs.Range("E10").Value = 9
$reboundsWs.Range("F10").Value = 7
$reboundsWs.Range("G10").Value = 3
$reboundsWs.Range("H10").Value = 1
$reboundsWs.Range("I10").Value = 0
$reboundsWs.Range("J10").Value = 1
$reboundsWs.Range("K10").Value = 7
$reboundsWs.Range("L10").Value = 4
$reboundsWs.Range("M10").Value = 0
$reboundsWs.Range("N10").Value = 0
$reboundsWs.Range("O10").Value = 0
$reboundsWs.Range("P10").Value = 1
$reboundsWs.Range("A11").Value = "'2025-11-07"
$reboundsWs.Range("B11").Value = "DEN"
$reboundsWs.Range("C11").Value = 2
$reboundsWs.Range("D11").Value = 7
$reboundsWs.Range("E11").Value = 2
$reboundsWs.Range("F11").Value = 3
$reboundsWs.Range("G11").Value = 3
$reboundsWs.Range("H11").Value = 0
$reboundsWs.Range("I11").Value = 3
$reboundsWs.Range("J11").Value = 4
$reboundsWs.Range("K11").Value = 0
$reboundsWs.Range("L11").Value = 3
$reboundsWs.Range("M11").Value = 6
$reboundsWs.Range("N11").Value = 0
$reboundsWs.Range("O11").Value = 7
$reboundsWs.Range("P11").Value = 6
$reboundsWs.Range("A12").Value = "'2025-11-09"
$reboundsWs.Range("B12").Value = "IND"
$reboundsWs.Range("C12").Value = 1
$reboundsWs.Range("D12").Value = 8
$reboundsWs.Range("E12").Value = 6
$reboundsWs.Range("F12").Value = 3
$reboundsWs.Range("G12").Value = 3
$reboundsWs.Range("H12").Value = 4
$reboundsWs.Range("I12").Value = 9
$reboundsWs.Range("J12").Value = 2
$reboundsWs.Range("K12").Value = 4
$reboundsWs.Range("L12").Value = 8
$reboundsWs.Range("M12").Value = 2
$reboundsWs.Range("N12").Value = 0
$reboundsWs.Range("O12").Value = 0
$reboundsWs.Range("P12").Value = 0

# --- Fill "3PM" data rows 2-12 ---
$pm3Ws.Range("A2").Value = "'2025-10-21"
$pm3Ws.Range("B2").Value = "LAL"
$pm3Ws.Range("C2").Value = 0
$pm3Ws.Range("D2").Value = 4
$pm3Ws.Range("E2").Value = 1
$pm3Ws.Range("F2").Value = 1
$pm3Ws.Range("G2").Value = 0
$pm3Ws.Range("H2").Value = 5
$pm3Ws.Range("I2").Value = 1
$pm3Ws.Range("J2").Value = 0
$pm3Ws.Range("K2").Value = 1
$pm3Ws.Range("L2").Value = 0
$pm3Ws.Range("M2").Value = 1
$pm3Ws.Range("N2").Value = 3
$pm3Ws.Range("O2").Value = 0
$pm3Ws.Range("P2").Value = 0
$pm3Ws.Range("A3").Value = "'2025-10-23"
$pm3Ws.Range("B3").Value = "DEN"
$pm3Ws.Range("C3").Value = 0
$pm3Ws.Range("D3").Value = 1
$pm3Ws.Range("E3").Value = 2
$pm3Ws.Range("F3").Value = 1
$pm3Ws.Range("G3").Value = 0
$pm3Ws.Range("H3").Value = 1
$pm3Ws.Range("I3").Value = 2
$pm3Ws.Range("J3").Value = 0
$pm3Ws.Range("K3").Value = 3
$pm3Ws.Range("L3").Value = 0
$pm3Ws.Range("M3").Value = 2
$pm3Ws.Range("N3").Value = 6
$pm3Ws.Range("O3").Value = 0
$pm3Ws.Range("P3").Value = 0
$pm3Ws.Range("A4").Value = "'2025-10-24"
$pm3Ws.Range("B4").Value = "POR"
$pm3Ws.Range("C4").Value = 0
$pm3Ws.Range("D4").Value = 1
$pm3Ws.Range("E4").Value = 1
$pm3Ws.Range("F4").Value = 1
$pm3Ws.Range("G4").Value = 1
$pm3Ws.Range("H4").Value = 0
$pm3Ws.Range("I4").Value = 0
$pm3Ws.Range("J4").Value = 1
$pm3Ws.Range("K4").Value = 0
$pm3Ws.Range("L4").Value = 1
$pm3Ws.Range("M4").Value = 3
$pm3Ws.Range("N4").Value = 7
$pm3Ws.Range("O4").Value = 0
$pm3Ws.Range("P4").Value = 0
$pm3Ws.Range("A5").Value = "'2025-10-27"
$pm3Ws.Range("B5").Value = "MEM"
$pm3Ws.Range("C5").Value = 0
$pm3Ws.Range("D5").Value = 0
$pm3Ws.Range("E5").Value = 5
$pm3Ws.Range("F5").Value = 0
$pm3Ws.Range("G5").Value = 5
$pm3Ws.Range("H5").Value = 0
$pm3Ws.Range("I5").Value = 2
$pm3Ws.Range("J5").Value = 1
$pm3Ws.Range("K5").Value = 0
$pm3Ws.Range("L5").Value = 0
$pm3Ws.Range("M5").Value = 2
$pm3Ws.Range("N5").Value = 4
$pm3Ws.Range("O5").Value = 0
$pm3Ws.Range("P5").Value = 0
$pm3Ws.Range("A6").Value = "'2025-10-28"
$pm3Ws.Range("B6").Value = "LAC"
$pm3Ws.Range("C6").Value = 0
$pm3Ws.Range("D6").Value = 1
$pm3Ws.Range("E6").Value = 2
$pm3Ws.Range("F6").Value = 1
$pm3Ws.Range("G6").Value = 3
$pm3Ws.Range("H6").Value = 0
$pm3Ws.Range("I6").Value = 3
$pm3Ws.Range("J6").Value = 0
$pm3Ws.Range("K6").Value = 0
$pm3Ws.Range("L6").Value = 4
$pm3Ws.Range("M6").Value = 1
$pm3Ws.Range("N6").Value = 2
$pm3Ws.Range("O6").Value = 0
$pm3Ws.Range("P6").Value = 0
$pm3Ws.Range("A7").Value = "'2025-10-30"
$pm3Ws.Range("B7").Value = "MIL"
$pm3Ws.Range("C7").Value = 0
$pm3Ws.Range("D7").Value = 2
$pm3Ws.Range("E7").Value = 1
$pm3Ws.Range("F7").Value = 1
$pm3Ws.Range("G7").Value = 0
$pm3Ws.Range("H7").Value = 2
$pm3Ws.Range("I7").Value = 1
$pm3Ws.Range("J7").Value = 0
$pm3Ws.Range("K7").Value = 1
$pm3Ws.Range("L7").Value = 0
$pm3Ws.Range("M7").Value = 1
$pm3Ws.Range("N7").Value = 4
$pm3Ws.Range("O7").Value = 0
$pm3Ws.Range("P7").Value = 0
$pm3Ws.Range("A8").Value = "'2025-11-01"
$pm3Ws.Range("B8").Value = "IND"
$pm3Ws.Range("C8").Value = 0
$pm3Ws.Range("D8").Value = 1
$pm3Ws.Range("E8").Value = 1
$pm3Ws.Range("F8").Value = 0
$pm3Ws.Range("G8").Value = 1
$pm3Ws.Range("H8").Value = 1
$pm3Ws.Range("I8").Value = 0
$pm3Ws.Range("J8").Value = 1
$pm3Ws.Range("K8").Value = 0
$pm3Ws.Range("L8").Value = 1
$pm3Ws.Range("M8").Value = 2
$pm3Ws.Range("N8").Value = 4
$pm3Ws.Range("O8").Value = 0
$pm3Ws.Range("P8").Value = 0
$pm3Ws.Range("A9").Value = "'2025-11-04"
$pm3Ws.Range("B9").Value = "PHX"
$pm3Ws.Range("C9").Value = 0
$pm3Ws.Range("D9").Value = 0
$pm3Ws.Range("E9").Value = 1
$pm3Ws.Range("F9").Value = 0
$pm3Ws.Range("G9").Value = 5
$pm3Ws.Range("H9").Value = 3
$pm3Ws.Range("I9").Value = 0
$pm3Ws.Range("J9").Value = 0
$pm3Ws.Range("K9").Value = 0
$pm3Ws.Range("L9").Value = 4
$pm3Ws.Range("M9").Value = 0
$pm3Ws.Range("N9").Value = 5
$pm3Ws.Range("O9").Value = 0
$pm3Ws.Range("P9").Value = 1
$pm3Ws.Range("A10").Value = "'2025-11-05"
$pm3Ws.Range("B10").Value = "SAC"
$pm3Ws.Range("C10").Value = 1
$pm3Ws.Range("D10").Value = 2
$pm3Ws.Range("E10").Value = 2
$pm3Ws.Range("F10").Value = 5
$pm3Ws.Range("G10").Value = 6
$pm3Ws.Range("H10").Value = 1
$pm3Ws.Range("I10").Value = 0
$pm3Ws.Range("J10").Value = 0
$pm3Ws.Range("K10").Value = 0
$pm3Ws.Range("L10").Value = 0
$pm3Ws.Range("M10").Value = 0
$pm3Ws.Range("N10").Value = 0
$pm3Ws.Range("O10").Value = 0
$pm3Ws.Range("P10").Value = 0
$pm3Ws.Range("A11").Value = "'2025-11-07"
$pm3Ws.Range("B11").Value = "DEN"
$pm3Ws.Range("C11").Value = 0
$pm3Ws.Range("D11").Value = 0
$pm3Ws.Range("E11").Value = 1
$pm3Ws.Range("F11").Value = 1
$pm3Ws.Range("G11").Value = 0
$pm3Ws.Range("H11").Value = 2
$pm3Ws.Range("I11").Value = 0
$pm3Ws.Range("J11").Value = 0
$pm3Ws.Range("K11").Value = 0
$pm3Ws.Range("L11").Value = 2
$pm3Ws.Range("M11").Value = 5
$pm3Ws.Range("N11").Value = 0
$pm3Ws.Range("O11").Value = 0
$pm3Ws.Range("P11").Value = 1
$pm3Ws.Range("A12").Value = "'2025-11-09"
$pm3Ws.Range("B12").Value = "IND"
$pm3Ws.Range("C12").Value = 0
$pm3Ws.Range("D12").Value = 0
$pm3Ws.Range("E12").Value = 1
$pm3Ws.Range("F12").Value = 0
$pm3Ws.Range("G12").Value = 3
$pm3Ws.Range("H12").Value = 0
$pm3Ws.Range("I12").Value = 0
$pm3Ws.Range("J12").Value = 1
$pm3Ws.Range("K12").Value = 4
$pm3Ws.Range("L12").Value = 2
$pm3Ws.Range("M12").Value = 1
$pm3Ws.Range("N12").Value = 0
$pm3Ws.Range("O12").Value = 0
$pm3Ws.Range("P12").Value = 0

# --- Fill "Avg Rebounds" data rows 2-15 ---
$avgReboundsWs.Range("A2").Value = "Jonathan Kuminga"
$avgReboundsWs.Range("B2").Value = 7.181818181818182
$avgReboundsWs.Range("A3").Value = "Draymond Green"
$avgReboundsWs.Range("B3").Value = 5.6
$avgReboundsWs.Range("A4").Value = "Jimmy Butler III"
$avgReboundsWs.Range("B4").Value = 5.5
$avgReboundsWs.Range("A5").Value = "Brandin Podziemski"
$avgReboundsWs.Range("B5").Value = 4.727272727272728
$avgReboundsWs.Range("A6").Value = "Al Horford"
$avgReboundsWs.Range("B6").Value = 4.428571428571429
$avgReboundsWs.Range("A7").Value = "Quinten Post"
$avgReboundsWs.Range("B7").Value = 4.181818181818182
$avgReboundsWs.Range("A8").Value = "Stephen Curry"
$avgReboundsWs.Range("B8").Value = 3.625
$avgReboundsWs.Range("A9").Value = "Moses Moody"
$avgReboundsWs.Range("B9").Value = 2.888888888888889
$avgReboundsWs.Range("A10").Value = "Will Richard"
$avgReboundsWs.Range("B10").Value = 2.555555555555555
$avgReboundsWs.Range("A11").Value = "Trayce Jackson-Davis"
$avgReboundsWs.Range("B11").Value = 2
$avgReboundsWs.Range("A12").Value = "Gary Payton II"
$avgReboundsWs.Range("B12").Value = 1.9
$avgReboundsWs.Range("A13").Value = "Pat Spencer"
$avgReboundsWs.Range("B13").Value = 1.714285714285714
$avgReboundsWs.Range("A14").Value = "Buddy Hield"
$avgReboundsWs.Range("B14").Value = 1.545454545454545
$avgReboundsWs.Range("A15").Value = "Gui Santos"
$avgReboundsWs.Range("B15").Value = 1.3

# --- Fill "Avg 3PM" data rows 2-15 ---
$avgPm3Ws.Range("A2").Value = "Stephen Curry"
$avgPm3Ws.Range("B2").Value = 4.375
$avgPm3Ws.Range("A3").Value = "Moses Moody"
$avgPm3Ws.Range("B3").Value = 2.666666666666667
$avgPm3Ws.Range("A4").Value = "Draymond Green"
$avgPm3Ws.Range("B4").Value = 1.8
$avgPm3Ws.Range("A5").Value = "Brandin Podziemski"
$avgPm3Ws.Range("B5").Value = 1.636363636363636
$avgPm3Ws.Range("A6").Value = "Buddy Hield"
$avgPm3Ws.Range("B6").Value = 1.363636363636364
$avgPm3Ws.Range("A7").Value = "Al Horford"
$avgPm3Ws.Range("B7").Value = 1.285714285714286
$avgPm3Ws.Range("A8").Value = "Quinten Post"
$avgPm3Ws.Range("B8").Value = 1.272727272727273
$avgPm3Ws.Range("A9").Value = "Will Richard"
$avgPm3Ws.Range("B9").Value = 1.222222222222222
$avgPm3Ws.Range("A10").Value = "Jonathan Kuminga"
$avgPm3Ws.Range("B10").Value = 1.090909090909091
$avgPm3Ws.Range("A11").Value = "Jimmy Butler III"
$avgPm3Ws.Range("B11").Value = 0.9
$avgPm3Ws.Range("A12").Value = "Gui Santos"
$avgPm3Ws.Range("B12").Value = 0.4
$avgPm3Ws.Range("A13").Value = "Pat Spencer"
$avgPm3Ws.Range("B13").Value = 0.2857142857142857
$avgPm3Ws.Range("A14").Value = "Gary Payton II"
$avgPm3Ws.Range("B14").Value = 0.1
$avgPm3Ws.Range("A15").Value = "Trayce Jackson-Davis"
$avgPm3Ws.Range("B15").Value = 0

Write-Output "done"
